# Apply cryptocurrency price/volume updates to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '49.841.71'
$ws.Range('D3').Value = '2.653.10'
$ws.Range('E3').Value = '  +0.14%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '''113.21'
$ws.Range('E5').Value = '  -0.66%  '
$ws.Range('D6').Value = '''327.99'
$ws.Range('E6').Value = '  +0.54%  '
$ws.Range('E7').Value = '  -1.07%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '''0.552'
$ws.Range('E9').Value = '  -1.38%  '
$ws.Range('D10').Value = '''39.93'
$ws.Range('E10').Value = '  -2.75%  '
$ws.Range('D11').Value = '''20.05'
$ws.Range('E11').Value = '  -0.36%  '
$ws.Range('E12').Value = '  -0.87%  '
$ws.Range('E13').Value = '  +1.87%  '
$ws.Range('D14').Value = '''7.60'
$ws.Range('E14').Value = '  +2.67%  '
$ws.Range('D15').Value = '3.067.73'
$ws.Range('E15').Value = '  +0.31%  '
$ws.Range('D16').Value = '2.649.03'
$ws.Range('E16').Value = '  +0.66%  '
$ws.Range('E17').Value = '  -1.62%  '
$ws.Range('D18').Value = '49.780.75'
$ws.Range('E18').Value = '  -0.31%  '
$ws.Range('D19').Value = '''13.41'
$ws.Range('E19').Value = '  +1.37%  '
$ws.Range('E20').Value = '  +0.84%  '
$ws.Range('E21').Value = '  -0.67%  '
$ws.Range('E22').Value = '  -0.72%  '
$ws.Range('D23').Value = '''269.46'
$ws.Range('E23').Value = '  -2.72%  '
$ws.Range('D24').Value = '''69.36'
$ws.Range('E24').Value = '  -4.20%  '
$ws.Range('E25').Value = '  -0.57%  '
$ws.Range('D26').Value = '''26.29'
$ws.Range('E26').Value = '  -2.40%  '
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('D28').Value = '''10.22'
$ws.Range('E28').Value = '  +1.89%  '
$ws.Range('E29').Value = '  -0.68%  '
$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D30').Value = '''0.138'
$ws.Range('E30').Value = '  -2.24%  '
$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D31').Value = '''34.99'
$ws.Range('E31').Value = '  -4.23%  '
$ws.Range('D32').Value = '''49.60'
$ws.Range('E32').Value = '  -1.15%  '
$ws.Range('D33').Value = '''5.52'
$ws.Range('E33').Value = '  +0.68%  '
$ws.Range('D34').Value = '''0.0821'
$ws.Range('E34').Value = '  +0.37%  '
$ws.Range('D35').Value = '''19.27'
$ws.Range('E35').Value = '  -0.99%  '
$ws.Range('E36').Value = '  -0.05%  '
$ws.Range('D37').Value = '''4.97'
$ws.Range('E37').Value = '  -1.50%  '
$ws.Range('E38').Value = '  -1.22%  '
$ws.Range('D39').Value = '''3.13'
$ws.Range('E39').Value = '  +0.62%  '
$ws.Range('B40').Value = 'Monero'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D40').Value = '''128.99'
$ws.Range('E40').Value = '  +3.68%  '
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').Value = '''23.72'
$ws.Range('E41').Value = '  +8.04%  '
$ws.Range('D42').Value = '''0.0345'
$ws.Range('E42').Value = '  +8.42%  '
$ws.Range('E43').Value = '  +2.53%  '
$ws.Range('E44').Value = '  -0.45%  '
$ws.Range('E45').Value = '  -0.17%  '
$ws.Range('D46').Value = '2.068.41'
$ws.Range('E47').Value = '  +7.62%  '
$ws.Range('E48').Value = '  -3.44%  '
$ws.Range('E49').Value = '  -2.13%  '
$ws.Range('D50').Value = '''5.29'
$ws.Range('E50').Value = '  -1.13%  '
$ws.Range('D51').Value = '''59.50'
$ws.Range('E51').Value = '  -0.49%  '
